$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.21"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.24"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.042"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05604"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.555"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.019"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8133"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8347"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1336"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06948"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03235"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02828"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09402"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001516"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005940"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006241"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.500"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.091"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3186"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1293"
$ws.Range("E21").Value = "20ProBitTokenPROB"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.741"
$ws.Range("E22").Value = "21MCDexMCB"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04675"
$ws.Range("E23").Value = "22CoinExTokenCET"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1370"
$ws.Range("E24").Value = "23ZBTokenZB"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001242"
$ws.Range("E25").Value = "24BitKanKAN"
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004286"
$ws.Range("E26").Value = "25HotbitTokenHTB"
$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009701"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001940"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03660"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006220"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1050"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002729"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008181"

Write-Host "Applied symbol list update"
